$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1999.7778
$ws.Range("I11").Value = 1999.7778
$ws.Range("K11").Value = 1999.7778
$ws.Range("M11").Value = -1859.7778
$ws.Range("H17").Value = 203194.58
$ws.Range("J17").Value = 209513.16
$ws.Range("L17").Value = 628539.48
$ws.Range("N17").Value = -628875.48
$ws.Range("H70").Value = 982.5
$ws.Range("I70").Value = 994.5
$ws.Range("J70").Value = 978.5
$ws.Range("K70").Value = 2983.5
$ws.Range("L70").Value = 2935.5
$ws.Range("M70").Value = -2713.5
$ws.Range("N70").Value = -3475.5
$ws.Range("H73").Value = 982.5
$ws.Range("I73").Value = 994.5
$ws.Range("J73").Value = 978.5
$ws.Range("K73").Value = 2983.5
$ws.Range("L73").Value = 2935.5
$ws.Range("M73").Value = -2047.5
$ws.Range("N73").Value = -4807.5
$ws.Range("H82").Value = 2566.7144
$ws.Range("I82").Value = 2566.7144
$ws.Range("K82").Value = 7700.1432
$ws.Range("M82").Value = -7294.1432
$ws.Range("H85").Value = 2566.7144
$ws.Range("I85").Value = 2566.7144
$ws.Range("K85").Value = 7700.1432
$ws.Range("M85").Value = -6296.1432
$ws.Range("H86").Value = 88892776
$ws.Range("I86").Value = 125003880
$ws.Range("J86").Value = 47622932
$ws.Range("K86").Value = 125003880
$ws.Range("L86").Value = 47622932
$ws.Range("M86").Value = -125002757
$ws.Range("N86").Value = -47625178
$ws.Range("H89").Value = 88892776
$ws.Range("I89").Value = 125003880
$ws.Range("J89").Value = 47622932
$ws.Range("K89").Value = 625019400
$ws.Range("L89").Value = 238114660
$ws.Range("M89").Value = -625013784
$ws.Range("N89").Value = -238125892
$ws.Range("H92").Value = 9377400
$ws.Range("I92").Value = 3909249.8
$ws.Range("K92").Value = 3909249.8
$ws.Range("M92").Value = -3908001.8
$ws.Range("H97").Value = 3553.4285
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 3553.4285
$ws.Range("K97").Value = 0
$ws.Range("L97").ClearContents()
$ws.Range("M97").Value = 10660.2855
$ws.Range("N97").Value = -11652.2855
$ws.Range("H106").Value = 47621570
$ws.Range("I106").Value = 50002424
$ws.Range("K106").Value = 50002424
$ws.Range("M106").Value = -50001793
$ws.Range("H116").Value = 5516.3335
$ws.Range("I116").Value = 5449.5
$ws.Range("J116").Value = 5549.75
$ws.Range("K116").Value = 5449.5
$ws.Range("L116").Value = 5549.75
$ws.Range("M116").Value = -2007.5
$ws.Range("N116").Value = -12433.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1426.4546
$ws.Range("J2").Value = 1613.5
$ws.Range("L2").Value = 1613.5
$ws.Range("N2").Value = -1839.5
$ws.Range("H32").Value = 153847.81
$ws.Range("I32").Value = 176709.38
$ws.Range("K32").Value = 176709.38
$ws.Range("M32").Value = -176422.38
$ws.Range("H61").Value = 1607165.9
$ws.Range("I61").Value = 54537.906
$ws.Range("K61").Value = 54537.906
$ws.Range("M61").Value = -54325.906
$ws.Range("H64").Value = 18000
$ws.Range("I64").Value = 18000
$ws.Range("K64").Value = 18000
$ws.Range("M64").Value = -17752
$ws.Range("H67").Value = 18000
$ws.Range("I67").Value = 18000
$ws.Range("K67").Value = 18000
$ws.Range("M67").Value = -17142
$ws.Range("H80").Value = 51990.5
$ws.Range("J80").Value = 51990.5
$ws.Range("L80").Value = 51990.5
$ws.Range("N80").Value = -53986.5
$ws.Range("H83").Value = 51990.5
$ws.Range("J83").Value = 51990.5
$ws.Range("L83").Value = 155971.5
$ws.Range("N83").Value = -165955.5
$ws.Range("H116").Value = 1426.4546
$ws.Range("J116").Value = 1613.5
$ws.Range("L116").Value = 1613.5
$ws.Range("N116").Value = -6201.5
$ws.Range("H132").Value = 1680.68
$ws.Range("I132").Value = 1198.5405
$ws.Range("J132").Value = 3052.923
$ws.Range("K132").Value = 3595.6215
$ws.Range("L132").Value = 9158.769
$ws.Range("M132").Value = -1065.6215
$ws.Range("N132").Value = -14218.769
$ws.Range("H136").Value = 1607165.9
$ws.Range("I136").Value = 54537.906
$ws.Range("K136").Value = 163613.718
$ws.Range("M136").Value = -161063.718

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1426.4546
$ws.Range("J3").Value = 1613.5
$ws.Range("L3").Value = 1613.5
$ws.Range("N3").Value = -1841.5
$ws.Range("H86").Value = 6366.826
$ws.Range("I86").Value = 4208.5713
$ws.Range("J86").Value = 7311.0625
$ws.Range("K86").Value = 4208.5713
$ws.Range("L86").Value = 7311.0625
$ws.Range("M86").Value = -3085.5713
$ws.Range("N86").Value = -9557.0625
$ws.Range("H89").Value = 6366.826
$ws.Range("I89").Value = 4208.5713
$ws.Range("J89").Value = 7311.0625
$ws.Range("K89").Value = 21042.8565
$ws.Range("L89").Value = 36555.3125
$ws.Range("M89").Value = -15426.8565
$ws.Range("N89").Value = -47787.3125
$ws.Range("H105").Value = 10434.667
$ws.Range("J105").Value = 11827.556
$ws.Range("L105").Value = 11827.556
$ws.Range("N105").Value = -15321.556

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 51169.4
$ws.Range("J44").Value = 63887
$ws.Range("L44").Value = 191661
$ws.Range("N44").Value = -192457
$ws.Range("H129").Value = 5563.7856
$ws.Range("I129").Value = 1362.1666
$ws.Range("J129").Value = 8715
$ws.Range("K129").Value = 4086.4998
$ws.Range("L129").Value = 26145
$ws.Range("M129").Value = 913.5001999999999
$ws.Range("N129").Value = -36145
$ws.Range("H133").Value = 4805.222
$ws.Range("I133").Value = 3912.75
$ws.Range("K133").Value = 11738.25
$ws.Range("M133").Value = -6678.25
$ws.Range("H137").Value = 1762.6
$ws.Range("I137").Value = 1703.3334
$ws.Range("K137").Value = 5110.0002
$ws.Range("M137").Value = -10.0002000000004
$ws.Range("H138").Value = 2906.611
$ws.Range("I138").Value = 2986.2354
$ws.Range("K138").Value = 8958.706200000001
$ws.Range("M138").Value = -3818.706200000001
$ws.Range("H139").Value = 4467167
$ws.Range("I139").Value = 5954365.5
$ws.Range("J139").Value = 5573
$ws.Range("K139").Value = 17863096.5
$ws.Range("L139").Value = 16719
$ws.Range("M139").Value = -17857956.5
$ws.Range("N139").Value = -26999
$ws.Range("H140").Value = 2520.8333
$ws.Range("I140").Value = 1943.1111
$ws.Range("J140").Value = 2954.125
$ws.Range("K140").Value = 5829.3333
$ws.Range("L140").Value = 8862.375
$ws.Range("M140").Value = -649.3333000000002
$ws.Range("N140").Value = -19222.375

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 89
$ws.Range("I2").Value = 84.42856999999999
$ws.Range("J2").Value = 97
$ws.Range("K2").Value = 84.42856999999999
$ws.Range("L2").Value = 97
$ws.Range("M2").Value = 28.57143000000001
$ws.Range("N2").Value = -323
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").ClearContents()
$ws.Range("N40").Value = 0
$ws.Range("H70").Value = 3083.5
$ws.Range("I70").Value = 3020.4443
$ws.Range("K70").Value = 3020.4443
$ws.Range("M70").Value = -2750.4443
$ws.Range("H73").Value = 3083.5
$ws.Range("I73").Value = 3020.4443
$ws.Range("K73").Value = 3020.4443
$ws.Range("M73").Value = -2084.4443
$ws.Range("H132").Value = 781900.2
$ws.Range("I132").Value = 7069.263
$ws.Range("J132").Value = 1833456.5
$ws.Range("K132").Value = 21207.789
$ws.Range("L132").Value = 5500369.5
$ws.Range("M132").Value = -18677.789
$ws.Range("N132").Value = -5505429.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11028.75
$ws.Range("J7").Value = 12762.454
$ws.Range("L7").Value = 12762.454
$ws.Range("N7").Value = -12986.454
$ws.Range("H22").Value = 3703.6667
$ws.Range("I22").Value = 1039.625
$ws.Range("J22").Value = 5035.6875
$ws.Range("K22").Value = 1039.625
$ws.Range("L22").Value = 5035.6875
$ws.Range("M22").Value = -744.625
$ws.Range("N22").Value = -5625.6875
$ws.Range("J25").Value = 11111
$ws.Range("L25").Value = 11111
$ws.Range("N25").Value = -11571
$ws.Range("H27").Value = 3703.6667
$ws.Range("I27").Value = 1039.625
$ws.Range("J27").Value = 5035.6875
$ws.Range("K27").Value = 1039.625
$ws.Range("L27").Value = 5035.6875
$ws.Range("M27").Value = -932.625
$ws.Range("N27").Value = -5249.6875
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("H55").Value = 396.14285
$ws.Range("J55").Value = 460.42105
$ws.Range("L55").Value = 460.42105
$ws.Range("N55").Value = -806.4210499999999
$ws.Range("H68").Value = 6027
$ws.Range("I68").Value = 6327.5713
$ws.Range("K68").Value = 6327.5713
$ws.Range("M68").Value = -5578.5713
$ws.Range("H71").Value = 6027
$ws.Range("I71").Value = 6327.5713
$ws.Range("K71").Value = 31637.8565
$ws.Range("M71").Value = -27893.8565
$ws.Range("H126").Value = 11028.75
$ws.Range("J126").Value = 12762.454
$ws.Range("L126").Value = 38287.362
$ws.Range("N126").Value = -43227.362
$ws.Range("H130").Value = 177777
$ws.Range("J130").Value = 177777
$ws.Range("L130").Value = 177777
$ws.Range("N130").Value = -187817

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3355.3333
$ws.Range("I126").Value = 3263.6155
$ws.Range("J126").Value = 3593.8
$ws.Range("K126").Value = 9790.8465
$ws.Range("L126").Value = 10781.4
$ws.Range("M126").Value = -7320.8465
$ws.Range("N126").Value = -15721.4
$ws.Range("H132").Value = 2882.7222
$ws.Range("I132").Value = 2579.6667
$ws.Range("K132").Value = 7739.000100000001
$ws.Range("M132").Value = -5209.000100000001
